$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.02211840812437771
$ws.Range("D2").Value = 0.2281605407615217
$ws.Range("E2").Value = 0.1655041902529852
$ws.Range("F2").Value = 1.044059152216271
$ws.Range("G2").Value = 0.5048982488050342
$ws.Range("H2").Value = 0.6485391751380547
$ws.Range("I2").Value = 0.4989997458086108
$ws.Range("J2").Value = 0.1632197117049579
$ws.Range("K2").Value = 0.8056589153039511
$ws.Range("M2").Value = 0.3129514844872716
$ws.Range("N2").Value = 1.210212824776477
$ws.Range("O2").Value = 2.271558306032517
$ws.Range("C3").Value = 0.01946547193403347
$ws.Range("D3").Value = 0.2258667498629023
$ws.Range("E3").Value = 0.1652913602818984
$ws.Range("F3").Value = 1.04579205144001
$ws.Range("G3").Value = 0.5054511691111259
$ws.Range("H3").Value = 0.6525337481336919
$ws.Range("I3").Value = 0.4990022720666083
$ws.Range("J3").Value = 0.1640221824738681
$ws.Range("K3").Value = 0.7091709022875534
$ws.Range("M3").Value = 0.2907857633141759
$ws.Range("N3").Value = 1.21121520712552
$ws.Range("O3").Value = 2.280829316501269
$ws.Range("C4").Value = 0.0178282052865697
$ws.Range("D4").Value = 0.224544235340062
$ws.Range("E4").Value = 0.1652350726342817
$ws.Range("F4").Value = 1.047448171944325
$ws.Range("G4").Value = 0.5061552958856268
$ws.Range("H4").Value = 0.6552835364291312
$ws.Range("I4").Value = 0.4993228547451842
$ws.Range("J4").Value = 0.1645973232094491
$ws.Range("K4").Value = 0.6497833931098285
$ws.Range("M4").Value = 0.2772331026052868
$ws.Range("N4").Value = 1.212242156522244
$ws.Range("O4").Value = 2.287907377829626
$ws.Range("C5").Value = 0.01715893373855693
$ws.Range("D5").Value = 0.2240269867784548
$ws.Range("E5").Value = 0.1652308724656493
$ws.Range("F5").Value = 1.048272002177853
$ws.Range("G5").Value = 0.5065338456468069
$ws.Range("H5").Value = 0.6564788520196529
$ws.Range("I5").Value = 0.4995337148107737
$ws.Range("J5").Value = 0.1648524289168982
$ws.Range("K5").Value = 0.6255482522211935
$ws.Range("M5").Value = 0.2717250874423485
$ws.Range("N5").Value = 1.212764416319423
$ws.Range("O5").Value = 2.291140096300751
$ws.Range("C6").Value = 0.01704767739249036
$ws.Range("D6").Value = 0.2239424105575551
$ws.Range("E6").Value = 0.1652313076512044
$ws.Range("F6").Value = 1.048417796280923
$ws.Range("G6").Value = 0.5066022344421768
$ws.Range("H6").Value = 0.6566818495865263
$ws.Range("I6").Value = 0.4995735729662911
$ws.Range("J6").Value = 0.1648960413121792
$ws.Range("K6").Value = 0.6215220109098141
$ws.Range("M6").Value = 0.2708113923741422
$ws.Range("N6").Value = 1.212857413467887
$ws.Range("O6").Value = 2.291697924082342
$ws.Range("C7").Value = 0.01781918761058421
$ws.Range("D7").Value = 0.2245371716145428
$ws.Range("E7").Value = 0.1652349400801434
$ws.Range("F7").Value = 1.04745867924958
$ws.Range("G7").Value = 0.5061600303075195
$ws.Range("H7").Value = 0.6552993541353516
$ws.Range("I7").Value = 0.4993253736743313
$ws.Range("J7").Value = 0.1646006797063819
$ws.Range("K7").Value = 0.6494566857692519
$ws.Range("M7").Value = 0.2771587590812317
$ws.Range("N7").Value = 1.212248779341849
$ws.Range("O7").Value = 2.287949565042908
$ws.Range("C8").Value = 0.02120542569902284
$ws.Range("D8").Value = 0.227351870226272
$ws.Range("E8").Value = 0.1654153850184734
$ws.Range("F8").Value = 1.044533793981522
$ws.Range("G8").Value = 0.505013168738941
$ws.Range("H8").Value = 0.6498548666622668
$ws.Range("I8").Value = 0.4989343985831098
$ws.Range("J8").Value = 0.1634793006254043
$ws.Range("K8").Value = 0.7724206954846693
$ws.Range("M8").Value = 0.3052971236859889
$ws.Range("N8").Value = 1.210473218838402
$ws.Range("O8").Value = 2.274467315300583
$ws.Range("C9").Value = 0.02777859969658891
$ws.Range("D9").Value = 0.2335493939344104
$ws.Range("E9").Value = 0.1663583158359216
$ws.Range("F9").Value = 1.04349498440758
$ws.Range("G9").Value = 0.5056615618921967
$ws.Range("H9").Value = 0.641533870500254
$ws.Range("I9").Value = 0.5006999565374528
$ws.Range("J9").Value = 0.161934089269117
$ws.Range("K9").Value = 1.012345229159052
$ws.Range("M9").Value = 0.3609141982558981
$ws.Range("N9").Value = 1.210244217481588
$ws.Range("O9").Value = 2.259028465020776
$ws.Range("C10").Value = 0.03256608436493025
$ws.Range("D10").Value = 0.2385119327819183
$ws.Range("E10").Value = 0.1674089124060103
$ws.Range("F10").Value = 1.045595162854788
$ws.Range("G10").Value = 0.5079111037872366
$ws.Range("H10").Value = 0.6368544529268974
$ws.Range("I10").Value = 0.5035431526221856
$ws.Range("J10").Value = 0.1611973355970804
$ws.Range("K10").Value = 1.18780550151132
$ws.Range("M10").Value = 0.4020257513762218
$ws.Range("N10").Value = 1.212044360880711
$ws.Range("O10").Value = 2.254401967138847
$ws.Range("C11").Value = 0.03473479761056808
$ws.Range("D11").Value = 0.2408575607213947
$ws.Range("E11").Value = 0.1679642986222127
$ws.Range("F11").Value = 1.047172468669999
$ws.Range("G11").Value = 0.509321028698011
$ws.Range("H11").Value = 0.6350366473030249
$ws.Range("I11").Value = 0.5051728607625705
$ws.Range("J11").Value = 0.1609486975327457
$ws.Range("K11").Value = 1.267435420670665
$ws.Range("M11").Value = 0.4207792752454509
$ws.Range("N11").Value = 1.21328779201832
$ws.Range("O11").Value = 2.253758138010767
$ws.Range("C12").Value = 0.03555469509247189
$ws.Range("D12").Value = 0.2417583810162256
$ws.Range("E12").Value = 0.1681857211539679
$ws.Range("F12").Value = 1.047859161448017
$ws.Range("G12").Value = 0.5099106244794882
$ws.Range("H12").Value = 0.634392960664357
$ws.Range("I12").Value = 0.505838374995129
$ws.Range("J12").Value = 0.1608669830643095
$ws.Range("K12").Value = 1.297560543230929
$ws.Range("M12").Value = 0.4278877983380909
$ws.Range("N12").Value = 1.213819438197163
$ws.Range("O12").Value = 2.253724543373181
$ws.Range("C13").Value = 0.03537817594660453
$ws.Range("D13").Value = 0.241563815147785
$ws.Range("E13").Value = 0.1681375402081819
$ws.Range("F13").Value = 1.047707294226655
$ws.Range("G13").Value = 0.509781166240046
$ws.Range("H13").Value = 0.6345296035273833
$ws.Range("I13").Value = 0.5056928926766986
$ws.Range("J13").Value = 0.1608840285154685
$ws.Range("K13").Value = 1.291073880466797
$ws.Range("M13").Value = 0.4263565479993332
$ws.Range("N13").Value = 1.213702239335191
$ws.Range("O13").Value = 2.253722427101479
$ws.Range("C14").Value = 0.03480227825939153
$ws.Range("D14").Value = 0.2409314201223367
$ws.Range("E14").Value = 0.1679822926989907
$ws.Range("F14").Value = 1.047227171742989
$ws.Range("G14").Value = 0.5093684184881937
$ws.Range("H14").Value = 0.6349827955843068
$ws.Range("I14").Value = 0.5052266433403574
$ws.Range("J14").Value = 0.1609417255387697
$ws.Range("K14").Value = 1.269914425206764
$ws.Range("M14").Value = 0.421363960649316
$ws.Range("N14").Value = 1.213330314114899
$ws.Range("O14").Value = 2.253751160400554
$ws.Range("C15").Value = 0.03444934787104614
$ws.Range("D15").Value = 0.2405456956799128
$ws.Range("E15").Value = 0.1678886450675812
$ws.Range("F15").Value = 1.046944724615543
$ws.Range("G15").Value = 0.5091228535842163
$ws.Range("H15").Value = 0.6352662062255803
$ws.Range("I15").Value = 0.5049473530115733
$ws.Range("J15").Value = 0.1609786865631051
$ws.Range("K15").Value = 1.256949812333346
$ws.Range("M15").Value = 0.4183067508793101
$ws.Range("N15").Value = 1.213110407590548
$ws.Range("O15").Value = 2.253796140253115
$ws.Range("C16").Value = 0.03242416715056606
$ws.Range("D16").Value = 0.2383604062051887
$ws.Range("E16").Value = 0.1673741728773699
$ws.Range("F16").Value = 1.045504594323965
$ws.Range("G16").Value = 0.5078267494059361
$ws.Range("H16").Value = 0.6369795042081137
$ws.Range("I16").Value = 0.5034434164649255
$ws.Range("J16").Value = 0.1612153254648589
$ws.Range("K16").Value = 1.182597549136858
$ws.Range("M16").Value = 0.4008011649500034
$ws.Range("N16").Value = 1.211971619978655
$ws.Range("O16").Value = 2.254473450563069
$ws.Range("C17").Value = 0.03117942046611688
$ws.Range("D17").Value = 0.2370423114737434
$ws.Range("E17").Value = 0.1670783804344644
$ws.Range("F17").Value = 1.044780389328032
$ws.Range("G17").Value = 0.5071307167337693
$ws.Range("H17").Value = 0.6381101622144598
$ws.Range("I17").Value = 0.5026069525622887
$ws.Range("J17").Value = 0.1613826533577409
$ws.Range("K17").Value = 1.136935317662903
$ws.Range("M17").Value = 0.3900749621917612
$ws.Range("N17").Value = 1.211381556644014
$ws.Range("O17").Value = 2.255263209002663
$ws.Range("C18").Value = 0.0304626170195661
$ws.Range("D18").Value = 0.2362924803586282
$ws.Range("E18").Value = 0.1669155415453716
$ws.Range("F18").Value = 1.044422386119358
$ws.Range("G18").Value = 0.5067667598611081
$ws.Range("H18").Value = 0.6387897504247348
$ws.Range("I18").Value = 0.5021574961524777
$ws.Range("J18").Value = 0.1614870392130143
$ws.Range("K18").Value = 1.110654048172137
$ws.Range("M18").Value = 0.3839104251482581
$ws.Range("N18").Value = 1.211082142581105
$ws.Range("O18").Value = 2.255854944667817
$ws.Range("C19").Value = 0.03021977332286951
$ws.Range("D19").Value = 0.2360400291585734
$ws.Range("E19").Value = 0.1668616605894186
$ws.Range("F19").Value = 1.044311227659513
$ws.Range("G19").Value = 0.5066497761472846
$ws.Range("H19").Value = 0.6390248742249867
$ws.Range("I19").Value = 0.5020107542898273
$ws.Range("J19").Value = 0.1615237811258652
$ws.Range("K19").Value = 1.101752706327829
$ws.Range("M19").Value = 0.3818240728921509
$ws.Range("N19").Value = 1.210987639898178
$ws.Range("O19").Value = 2.256078905191885
$ws.Range("C20").Value = 0.03131201492271884
$ws.Range("D20").Value = 0.2371817663455005
$ws.Range("E20").Value = 0.1671091134499996
$ws.Range("F20").Value = 1.044851423697231
$ws.Range("G20").Value = 0.5072010444285695
$ws.Range("H20").Value = 0.6379867732657658
$ws.Range("I20").Value = 0.5026927192265447
$ws.Range("J20").Value = 0.1613639982289641
$ws.Range("K20").Value = 1.141797972138079
$ws.Range("M20").Value = 0.3912162819092515
$ws.Range("N20").Value = 1.211440234569807
$ws.Range("O20").Value = 2.25516490733105
$ws.Range("C21").Value = 0.03497147019579927
$ws.Range("D21").Value = 0.2411168291443175
$ws.Range("E21").Value = 0.1680275913737326
$ws.Range("F21").Value = 1.047365769279779
$ws.Range("G21").Value = 0.5094881404015439
$ws.Range("H21").Value = 0.6348484697962959
$ws.Range("I21").Value = 0.5053622789532142
$ws.Range("J21").Value = 0.1609244409348065
$ws.Range("K21").Value = 1.276130272999296
$ws.Range("M21").Value = 0.4228302190581559
$ws.Range("N21").Value = 1.213437909922604
$ws.Range("O21").Value = 2.253737014529577
$ws.Range("C22").Value = 0.03735526799206923
$ws.Range("D22").Value = 0.2437619138745646
$ws.Range("E22").Value = 0.1686926126085773
$ws.Range("F22").Value = 1.04953012352199
$ws.Range("G22").Value = 0.5113075244324392
$ws.Range("H22").Value = 0.6330578052254339
$ws.Range("I22").Value = 0.5073889879067508
$ws.Range("J22").Value = 0.1607096703109576
$ws.Range("K22").Value = 1.363754245067128
$ws.Range("M22").Value = 0.443532171472917
$ws.Range("N22").Value = 1.215097694300326
$ws.Range("O22").Value = 2.254029133163641
$ws.Range("C23").Value = 0.0360837211536591
$ws.Range("D23").Value = 0.2423435067681794
$ws.Range("E23").Value = 0.1683317640527022
$ws.Range("F23").Value = 1.048327298157574
$ws.Range("G23").Value = 0.5103067496114733
$ws.Range("H23").Value = 0.6339896994065555
$ws.Range("I23").Value = 0.5062814873507548
$ws.Range("J23").Value = 0.1608176636504766
$ws.Range("K23").Value = 1.317003888733666
$ws.Range("M23").Value = 0.4324796035196883
$ws.Range("N23").Value = 1.214179513203234
$ws.Range("O23").Value = 2.253761057508711
$ws.Range("C24").Value = 0.03125207264729113
$ws.Range("D24").Value = 0.2371186939989656
$ws.Range("E24").Value = 0.1670951965742589
$ws.Range("F24").Value = 1.044819127295227
$ws.Range("G24").Value = 0.5071691365198774
$ws.Range("H24").Value = 0.6380424653362979
$ws.Range("I24").Value = 0.5026538461926435
$ws.Range("J24").Value = 0.1613724067117133
$ws.Range("K24").Value = 1.139599657613701
$ws.Range("M24").Value = 0.3907002846895224
$ws.Range("N24").Value = 1.211413582218938
$ws.Range("O24").Value = 2.255208920623943
$ws.Range("C25").Value = 0.02600766573807789
$ws.Range("D25").Value = 0.231800643770157
$ws.Range("E25").Value = 0.1660402965494292
$ws.Range("F25").Value = 1.043273167878134
$ws.Range("G25").Value = 0.505175257557255
$ws.Range("H25").Value = 0.643532945021704
$ws.Range("I25").Value = 0.4999509888689602
$ws.Range("J25").Value = 0.1622821217780164
$ws.Range("K25").Value = 0.9475764465604186
$ws.Range("M25").Value = 0.3458232157546774
$ws.Range("N25").Value = 1.209959460834455
$ws.Range("O25").Value = 2.262026368903577
